$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match rows (round 18) to append starting at row 170
$newRows = @(
    @{ Row=170; Idx=168; Home="Cagliari"; Away="Inter";         HS=0; AS=3; HXg=0.31; AXg=3.52; HXgP=0.4;  AXgP=3.33; HP=0; AP=2; HD=0.09; AD=0.19; TD=0.29; GD=1 },
    @{ Row=171; Idx=169; Home="Empoli";   Away="Genoa";         HS=1; AS=2; HXg=1.34; AXg=1.77; HXgP=1.42; AXgP=1.77; HP=1; AP=1; HD=0.08; AD=0;    TD=0.09; GD=1 },
    @{ Row=172; Idx=170; Home="Lazio";    Away="Atalanta";      HS=1; AS=1; HXg=0.62; AXg=2.88; HXgP=0.53; AXgP=2.28; HP=0; AP=1; HD=0.09; AD=0.6;  TD=0.7;  GD=1 },
    @{ Row=173; Idx=171; Home="Parma";    Away="Monza";         HS=2; AS=1; HXg=1.36; AXg=2.63; HXgP=1.63; AXgP=1.08; HP=1; AP=0; HD=0.27; AD=1.55; TD=1.82; GD=2 },
    @{ Row=174; Idx=172; Home="Juventus"; Away="Fiorentina";    HS=2; AS=2; HXg=1.49; AXg=1;    HXgP=1.35; AXgP=1.43; HP=0; AP=0; HD=0.14; AD=0.43; TD=0.57; GD=4 },
    @{ Row=175; Idx=173; Home="Milan";    Away="Roma";          HS=1; AS=1; HXg=1.88; AXg=1.26; HXgP=2.12; AXgP=1.11; HP=0; AP=0; HD=0.24; AD=0.15; TD=0.39; GD=2 },
    @{ Row=176; Idx=174; Home="Napoli";   Away="Venezia";       HS=1; AS=0; HXg=1.97; AXg=0.21; HXgP=2.49; AXgP=0.26; HP=1; AP=0; HD=0.52; AD=0.05; TD=0.57; GD=0 },
    @{ Row=177; Idx=175; Home="Udinese";  Away="Torino";        HS=2; AS=2; HXg=0.68; AXg=0.6;  HXgP=0.84; AXgP=0.82; HP=0; AP=0; HD=0.16; AD=0.22; TD=0.38; GD=4 },
    @{ Row=178; Idx=176; Home="Bologna";  Away="Hellas Verona"; HS=2; AS=3; HXg=2.32; AXg=0.84; HXgP=2.59; AXgP=0.8100000000000001; HP=0; AP=0; HD=0.27; AD=0.03; TD=0.3;  GD=5 },
    @{ Row=179; Idx=177; Home="Como";     Away="Lecce";         HS=2; AS=0; HXg=2.03; AXg=0.43; HXgP=2.01; AXgP=0.57; HP=1; AP=0; HD=0.02; AD=0.14; TD=0.16; GD=1 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Idx
    $ws.Cells.Item($row, 2).Value = $r.Home
    $ws.Cells.Item($row, 3).Value = $r.Away
    $ws.Cells.Item($row, 4).Value = $r.HS
    $ws.Cells.Item($row, 5).Value = $r.AS
    $ws.Cells.Item($row, 6).Value = $r.HXg
    $ws.Cells.Item($row, 7).Value = $r.AXg
    $ws.Cells.Item($row, 8).Value = $r.HXgP
    $ws.Cells.Item($row, 9).Value = $r.AXgP
    $ws.Cells.Item($row, 10).Value = $r.HP
    $ws.Cells.Item($row, 11).Value = $r.AP
    $ws.Cells.Item($row, 12).Value = $r.HD
    $ws.Cells.Item($row, 13).Value = $r.AD
    $ws.Cells.Item($row, 14).Value = $r.TD
    $ws.Cells.Item($row, 15).Value = $r.GD
}

# Match formatting of column A (bold, bordered, centered) used for the existing index column
$ws.Range("A169").Copy() | Out-Null
$ws.Range("A170:A179").PasteSpecial(-4122) | Out-Null
